$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.08112184783705709
$ws.Range("AQ2").Value = 1388.7
$ws.Range("AZ2").Value = 740.0999999999999
$ws.Range("BC2").Value = 164.3
$ws.Range("B3").Value = 0.01842065120159805
$ws.Range("C3").Value = -0.0238762625218272
$ws.Range("AQ3").Value = 1428.1
$ws.Range("AZ3").Value = 746.7
$ws.Range("BC3").Value = 164.9
$ws.Range("B4").Value = 0.4273326313133015
$ws.Range("C4").Value = 0.06409033606477599
$ws.Range("AQ4").Value = 1431.3
$ws.Range("AZ4").Value = 757.1999999999999
$ws.Range("BC4").Value = 164.5999999999999
$ws.Range("B5").Value = 0.1652983376755887
$ws.Range("C5").Value = 0.162628328928369
$ws.Range("AQ5").Value = 1434.8
$ws.Range("AZ5").Value = 760.6999999999999
$ws.Range("BC5").Value = 165.6
$ws.Range("B6").Value = -0.4280824757189849
$ws.Range("C6").Value = 0.04574228611787742
$ws.Range("AQ6").Value = 1440.2
$ws.Range("AZ6").Value = 755.8
$ws.Range("BC6").Value = 166
$ws.Range("B7").Value = 0.211407231042727
$ws.Range("C7").Value = 0.09398893107815967
$ws.Range("AQ7").Value = 1503.4
$ws.Range("AZ7").Value = 766.7
$ws.Range("BC7").Value = 167.3000000000001
$ws.Range("B8").Value = 0.6139640220240616
$ws.Range("C8").Value = 0.1406467787558497
$ws.Range("AQ8").Value = 1508.6
$ws.Range("AZ8").Value = 782.7
$ws.Range("BC8").Value = 167.7
$ws.Range("B9").Value = 0.4001201119220801
$ws.Range("C9").Value = 0.1993522223174726
$ws.Range("AQ9").Value = 1513.9
$ws.Range("AZ9").Value = 791
$ws.Range("BC9").Value = 168.7000000000001
$ws.Range("B10").Value = 0.4816080942708375
$ws.Range("C10").Value = 0.4267748648149282
$ws.Range("AQ10").Value = 1521.9
$ws.Range("AZ10").Value = 788.9
$ws.Range("BC10").Value = 169.5
$ws.Range("B11").Value = 0.08506766165855166
$ws.Range("C11").Value = 0.3951899724688844
$ws.Range("AQ11").Value = 1574.4
$ws.Range("AZ11").Value = 810.5
$ws.Range("BC11").Value = 186.4
$ws.Range("B12").Value = 14.28977758476561
$ws.Range("C12").Value = 3.814143363154272
$ws.Range("AE12").Value = 4134.6
$ws.Range("AQ12").Value = 3741.533333333334
$ws.Range("AZ12").Value = 1061.866666666667
$ws.Range("BA12").Value = 249.0666666666667
$ws.Range("BC12").Value = 393.0666666666668
$ws.Range("B13").Value = 3.607238917645103
$ws.Range("C13").Value = 4.615923064585027
$ws.Range("AO13").Value = 695.2
$ws.Range("AQ13").Value = 2615.2
$ws.Range("AZ13").Value = 865.5
$ws.Range("BA13").Value = 186
$ws.Range("BC13").Value = 181.8
$ws.Range("B14").Value = -2.144250383721548
$ws.Range("C14").Value = 3.95945844508693
$ws.Range("O14").Value = 298.4243775484318
$ws.Range("AC14").Value = 874.878
$ws.Range("AD14").Value = 3981.231964941214
$ws.Range("AE14").Value = 2200.735436095733
$ws.Range("AI14").Value = 5205.9886129062
$ws.Range("AJ14").Value = 309.2395471855706
$ws.Range("AO14").Value = 738.2
$ws.Range("AQ14").Value = 1932.605981095733
$ws.Range("AZ14").Value = 867.6637499999999
$ws.Range("BA14").Value = 140
$ws.Range("BC14").Value = 268.129455
$ws.Range("BG14").Value = 1909.190833092659
$ws.Range("BH14").Value = 92.02399613592002
$ws.Range("B15").Value = 3.383617228811452
$ws.Range("C15").Value = 4.784095836875155
$ws.Range("K15").Value = 1504.35145528078
$ws.Range("O15").Value = 341.1926948948232
$ws.Range("R15").Value = 2341.26256
$ws.Range("AB15").Value = 1585.922904250773
$ws.Range("AC15").Value = -6.476
$ws.Range("AD15").Value = 3709.608928678839
$ws.Range("AE15").Value = 3299.545529866908
$ws.Range("AI15").Value = 5266.950157518768
$ws.Range("AJ15").Value = 320.9813082499604
$ws.Range("AL15").Value = 664
$ws.Range("AN15").Value = 1399.852444352445
$ws.Range("AO15").Value = 43
$ws.Range("AQ15").Value = 3057.963251229407
$ws.Range("AZ15").Value = 869.8329093749999
$ws.Range("BA15").Value = -46
$ws.Range("BC15").Value = 241.5822786375001
$ws.Range("BG15").Value = 1940.854597891685
$ws.Range("BH15").Value = 93.55020615065929
$ws.Range("B16").Value = -5.508113222542727
$ws.Range("C16").Value = -0.1653768649519298
$ws.Range("K16").Value = 1511.252602267454
$ws.Range("O16").Value = 314.074370538223
$ws.Range("R16").Value = 2347.1157164
$ws.Range("AB16").Value = 1613.876822716431
$ws.Range("AC16").Value = 344.8
$ws.Range("AD16").Value = 3488.775830340336
$ws.Range("AE16").Value = 2053.740556383055
$ws.Range("AI16").Value = 5327.110323740309
$ws.Range("AJ16").Value = 332.6844424570025
$ws.Range("AN16").Value = 1427.806362818103
$ws.Range("AO16").Value = 257.8
$ws.Range("AQ16").Value = 1827.704322048961
$ws.Range("AZ16").Value = 872.0074916484373
$ws.Range("BA16").Value = 87
$ws.Range("BC16").Value = 226.0362343340938
$ws.Range("BG16").Value = 1971.716984299685
$ws.Range("BH16").Value = 95.03778930805085
$ws.Range("B17").Value = 2.048668124438896
$ws.Range("C17").Value = -0.5550195632534816
$ws.Range("K17").Value = 1590.736803023271
$ws.Range("O17").Value = 315.107190577634
$ws.Range("AB17").Value = 1636.87652348763
$ws.Range("AC17").Value = 724
$ws.Range("AD17").Value = 3333.663631241299
$ws.Range("AE17").Value = 1868.162483675872
$ws.Range("AI17").Value = 5389.395074533414
$ws.Range("AJ17").Value = 344.4899827509665
$ws.Range("AN17").Value = 1450.806063589301
$ws.Range("AQ17").Value = 1641.216068170108
$ws.Range("AZ17").Value = 876.3675291066794
$ws.Range("BA17").Value = 60
$ws.Range("BC17").Value = 226.9464155057643
$ws.Range("BG17").Value = 2004.703955279249
$ws.Range("BH17").Value = 96.62777855236428
$ws.Range("B18").Value = -1.058292142645675
$ws.Range("C18").Value = -0.2835300029845134
$ws.Range("K18").Value = 1602.397960524671
$ws.Range("O18").Value = 310.5969948902631
$ws.Range("AB18").Value = 1673.784384760761
$ws.Range("AC18").Value = 152
$ws.Range("AD18").Value = 3425.306801055291
$ws.Range("AE18").Value = 1864.902082462469
$ws.Range("AI18").Value = 5460.056646033359
$ws.Range("AJ18").Value = 388.860401143952
$ws.Range("AN18").Value = 1487.713924862433
$ws.Range("AQ18").Value = 1636.583568840412
$ws.Range("AZ18").Value = 882.9402855749795
$ws.Range("BC18").Value = 228.3185136220575
$ws.Range("BG18").Value = 2026.797771152849
$ws.Range("BH18").Value = 97.69271202645102
$ws.Range("B19").Value = -7.398079846028857
$ws.Range("C19").Value = -2.97895427169459
$ws.Range("K19").Value = 1574.879991543916
$ws.Range("O19").Value = 306.2938785049959
$ws.Range("AB19").Value = 1753.631582933507
$ws.Range("AD19").Value = 3490.923532880751
$ws.Range("AE19").Value = 1882.821171489655
$ws.Range("AN19").Value = 1524.621786135565
$ws.Range("AQ19").Value = 1652.659472731377
$ws.Range("AZ19").Value = 891.7696884307293
$ws.Range("BC19").Value = 230.1616987582779
$ws.Range("B20").Value = -2.312473296450245
$ws.Range("C20").Value = -2.18004429017147
$ws.Range("K20").Value = 1588.32908484966
$ws.Range("O20").Value = 308.01210563686
$ws.Range("AB20").Value = 1790.539444206639
$ws.Range("AD20").Value = 3539.116350869613
$ws.Range("AE20").Value = 1922.169636453873
$ws.Range("AN20").Value = 1561.529647408696
$ws.Range("AQ20").Value = 1651.524008138567
$ws.Range("AZ20").Value = 901.1633819594565
$ws.Range("BC20").Value = 270.6456283153065
$ws.Range("B21").Value = 1.409551007941173
$ws.Range("C21").Value = -2.3398235692959
$ws.Range("O21").Value = 309.9208485534183
$ws.Range("AB21").Value = 1825.097655294567
$ws.Range("AD21").Value = 3561.473045106598
$ws.Range("AE21").Value = 1935.531548589725
$ws.Range("AN21").Value = 1598.437508681828
$ws.Range("AQ21").Value = 1626.556319089618
$ws.Range("AZ21").Value = 909.6288449023414
$ws.Range("BC21").Value = 308.9752295001069
